# Daily attendance processing - normalize "Recorded By" (column G) entries.
# Entries that begin with "System, " have the leading "System" marker moved
# to the end of the comma-separated list (and any trailing lowercase
# "system" token promoted to "System" in the process), e.g.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com, system" -> "system, backup@backdoor.com, System"
# This is achieved generically by reversing the order of the comma-separated
# tokens in any "Recorded By" cell that currently starts with "System, ".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "System, *") {
        $parts = $val -split ", "
        $newVal = $parts[($parts.Count - 1)..0] -join ", "
        $cell.Value2 = $newVal
    }
}
